$d = $word.ActiveDocument

function Merge-IdRun([string]$token) {
    $target = "<id>$token</id>"

    $found = $d.Content.Find.Execute($target, $false, $false, $false, $false,
                                      $false, $true, 1, $false, "", 0)

    $range = $d.Content
    $range.Find.Execute($target, $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 0)
    if (-not $range.Find.Found) {
        throw "Could not find paragraph text '$target'"
    }

    $start = $range.Start
    $end = $range.End

    $mergedRange = $d.Range($start, $end)
    $mergedRange.Delete()

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;id&gt;' + $token + '&lt;/id&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $mergedRange.InsertXML($xml)
}

Merge-IdRun "p020r_1"
Merge-IdRun "p020r_2"
